$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 1 (A1:E1) values from 10 to 15
$ws.Range("A1:E1").Value = 15
